# Applies the crypto price/volume refresh described in the commit diff.
# Values in column D that look like plain decimal numbers are written with a
# leading apostrophe (forcing Excel to store them as text, preserving exact
# formatting such as trailing zeros) and then the style is reset to "Normal"
# so no stray number-format style gets attached to the cell.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Cell = "D2"; Value = "68.187.34"; ForceText = $false }
    @{ Cell = "E2"; Value = "  +1.06%  "; ForceText = $false }
    @{ Cell = "D3"; Value = "3.732.97"; ForceText = $false }
    @{ Cell = "E3"; Value = "  +0.38%  "; ForceText = $false }
    @{ Cell = "E4"; Value = "  +0.09%  "; ForceText = $false }
    @{ Cell = "D5"; Value = "'592.20"; ForceText = $true }
    @{ Cell = "E5"; Value = "  +0.51%  "; ForceText = $false }
    @{ Cell = "D6"; Value = "'166.84"; ForceText = $true }
    @{ Cell = "E6"; Value = "  +1.21%  "; ForceText = $false }
    @{ Cell = "D7"; Value = "3.732.25"; ForceText = $false }
    @{ Cell = "E7"; Value = "  +0.42%  "; ForceText = $false }
    @{ Cell = "E8"; Value = "  -0.04%  "; ForceText = $false }
    @{ Cell = "E9"; Value = "  +0.79%  "; ForceText = $false }
    @{ Cell = "E10"; Value = "  +0.96%  "; ForceText = $false }
    @{ Cell = "E11"; Value = "  +0.26%  "; ForceText = $false }
    @{ Cell = "D12"; Value = "'0.448"; ForceText = $true }
    @{ Cell = "E12"; Value = "  +0.14%  "; ForceText = $false }
    @{ Cell = "E13"; Value = "  -0.46%  "; ForceText = $false }
    @{ Cell = "D14"; Value = "'36.11"; ForceText = $true }
    @{ Cell = "E14"; Value = "  +0.93%  "; ForceText = $false }
    @{ Cell = "D15"; Value = "4.361.00"; ForceText = $false }
    @{ Cell = "E15"; Value = "  +0.43%  "; ForceText = $false }
    @{ Cell = "D16"; Value = "3.761.52"; ForceText = $false }
    @{ Cell = "E16"; Value = "  +1.16%  "; ForceText = $false }
    @{ Cell = "D17"; Value = "68.161.47"; ForceText = $false }
    @{ Cell = "E17"; Value = "  +1.13%  "; ForceText = $false }
    @{ Cell = "D18"; Value = "'17.82"; ForceText = $true }
    @{ Cell = "E18"; Value = "  -2.34%  "; ForceText = $false }
    @{ Cell = "B19"; Value = "TRON"; ForceText = $false }
    @{ Cell = "C19"; Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"; ForceText = $false }
    @{ Cell = "D19"; Value = "'0.112"; ForceText = $true }
    @{ Cell = "E19"; Value = "  +0.82%  "; ForceText = $false }
    @{ Cell = "B20"; Value = "Polkadot"; ForceText = $false }
    @{ Cell = "C20"; Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"; ForceText = $false }
    @{ Cell = "D20"; Value = "'6.99"; ForceText = $true }
    @{ Cell = "E20"; Value = "  -0.18%  "; ForceText = $false }
    @{ Cell = "D21"; Value = "'10.67"; ForceText = $true }
    @{ Cell = "E21"; Value = "  +0.86%  "; ForceText = $false }
    @{ Cell = "D22"; Value = "'465.37"; ForceText = $true }
    @{ Cell = "E22"; Value = "  +0.36%  "; ForceText = $false }
    @{ Cell = "D23"; Value = "'0.694"; ForceText = $true }
    @{ Cell = "E23"; Value = "  -0.40%  "; ForceText = $false }
    @{ Cell = "E24"; Value = "  +9.31%  "; ForceText = $false }
    @{ Cell = "D25"; Value = "'83.76"; ForceText = $true }
    @{ Cell = "E25"; Value = "  +1.82%  "; ForceText = $false }
    @{ Cell = "D26"; Value = "'2.17"; ForceText = $true }
    @{ Cell = "E26"; Value = "  +0.37%  "; ForceText = $false }
    @{ Cell = "D27"; Value = "'11.85"; ForceText = $true }
    @{ Cell = "E27"; Value = "  -0.45%  "; ForceText = $false }
    @{ Cell = "D28"; Value = "'10.03"; ForceText = $true }
    @{ Cell = "E28"; Value = "  -1.25%  "; ForceText = $false }
    @{ Cell = "E29"; Value = "  -0.03%  "; ForceText = $false }
    @{ Cell = "E30"; Value = "  +0.06%  "; ForceText = $false }
    @{ Cell = "D31"; Value = "'7.27"; ForceText = $true }
    @{ Cell = "E31"; Value = "  -0.29%  "; ForceText = $false }
    @{ Cell = "D32"; Value = "'29.79"; ForceText = $true }
    @{ Cell = "E32"; Value = "  +0.78%  "; ForceText = $false }
    @{ Cell = "E33"; Value = "  -2.04%  "; ForceText = $false }
    @{ Cell = "D34"; Value = "'9.15"; ForceText = $true }
    @{ Cell = "E34"; Value = "  +2.09%  "; ForceText = $false }
    @{ Cell = "D36"; Value = "3.686.90"; ForceText = $false }
    @{ Cell = "E36"; Value = "  +0.43%  "; ForceText = $false }
    @{ Cell = "D37"; Value = "'0.100"; ForceText = $true }
    @{ Cell = "E37"; Value = "  -0.33%  "; ForceText = $false }
    @{ Cell = "D38"; Value = "'3.42"; ForceText = $true }
    @{ Cell = "E38"; Value = "  +0.62%  "; ForceText = $false }
    @{ Cell = "E39"; Value = "  +1.91%  "; ForceText = $false }
    @{ Cell = "D40"; Value = "'0.994"; ForceText = $true }
    @{ Cell = "E40"; Value = "  +0.67%  "; ForceText = $false }
    @{ Cell = "D41"; Value = "'5.77"; ForceText = $true }
    @{ Cell = "E41"; Value = "  +1.02%  "; ForceText = $false }
    @{ Cell = "E42"; Value = "  +0.16%  "; ForceText = $false }
    @{ Cell = "D44"; Value = "'43.78"; ForceText = $true }
    @{ Cell = "E44"; Value = "  +16.03%  "; ForceText = $false }
    @{ Cell = "D45"; Value = "'0.299"; ForceText = $true }
    @{ Cell = "E45"; Value = "  -1.33%  "; ForceText = $false }
    @{ Cell = "D46"; Value = "'46.59"; ForceText = $true }
    @{ Cell = "E46"; Value = "  +3.10%  "; ForceText = $false }
    @{ Cell = "D47"; Value = "'1.90"; ForceText = $true }
    @{ Cell = "E47"; Value = "  -0.21%  "; ForceText = $false }
    @{ Cell = "D48"; Value = "'8.42"; ForceText = $true }
    @{ Cell = "E48"; Value = "  -0.65%  "; ForceText = $false }
    @{ Cell = "D49"; Value = "'387.86"; ForceText = $true }
    @{ Cell = "E49"; Value = "  -1.09%  "; ForceText = $false }
    @{ Cell = "D50"; Value = "'143.64"; ForceText = $true }
    @{ Cell = "E50"; Value = "  +0.23%  "; ForceText = $false }
    @{ Cell = "D51"; Value = "2.741.46"; ForceText = $false }
    @{ Cell = "E51"; Value = "  +2.57%  "; ForceText = $false }
)

foreach ($u in $updates) {
    $range = $ws.Range($u.Cell)
    $range.Value = $u.Value
    if ($u.ForceText) {
        $range.Style = "Normal"
    }
}
